# Edit script: add "clean_lifestyle" column (N) to Species_information worksheet,
# derived per-row from the existing "Pathogen"/"lifestyle" columns, and remove a
# stray Pathogen value in row 60 that the author cleaned up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new table column "clean_lifestyle" -----------------------------
$table = $ws.ListObjects.Item("Table1")
$table.ListColumns.Add() | Out-Null

# Header cell (bold, matches the other header cells)
$ws.Range("N1").Value = "clean_lifestyle"
$ws.Range("N1").Font.Bold = $true

# --- Per-row values for the new column --------------------------------------
$rows = @(
    @{Row=2; Value="Saprotroph"; Style2=$false},
    @{Row=3; Value="Plant pathogen"; Style2=$false},
    @{Row=4; Value="Plant pathogen"; Style2=$false},
    @{Row=5; Value="Insect pathogen"; Style2=$false},
    @{Row=6; Value="Plant pathogen"; Style2=$false},
    @{Row=7; Value="Human pathogen"; Style2=$false},
    @{Row=8; Value="Plant pathogen"; Style2=$false},
    @{Row=9; Value="Plant pathogen"; Style2=$false},
    @{Row=10; Value="Endophyte"; Style2=$false},
    @{Row=11; Value="Insect pathogen"; Style2=$false},
    @{Row=12; Value="Plant pathogen"; Style2=$false},
    @{Row=13; Value="Plant pathogen"; Style2=$false},
    @{Row=14; Value="Plant pathogen"; Style2=$false},
    @{Row=15; Value="Plant pathogen"; Style2=$false},
    @{Row=16; Value="Plant pathogen"; Style2=$false},
    @{Row=17; Value="Human pathogen"; Style2=$false},
    @{Row=18; Value="Human pathogen"; Style2=$false},
    @{Row=19; Value="Plant pathogen"; Style2=$false},
    @{Row=20; Value="Endophyte"; Style2=$false},
    @{Row=21; Value="Saprotroph"; Style2=$false},
    @{Row=22; Value="Plant pathogen"; Style2=$false},
    @{Row=23; Value="Plant pathogen"; Style2=$false},
    @{Row=24; Value="Plant pathogen"; Style2=$false},
    @{Row=25; Value="Insect pathogen"; Style2=$false},
    @{Row=26; Value="Plant pathogen"; Style2=$false},
    @{Row=27; Value="Human pathogen"; Style2=$false},
    @{Row=28; Value="Human pathogen"; Style2=$false},
    @{Row=29; Value="Plant pathogen"; Style2=$false},
    @{Row=30; Value="Plant pathogen"; Style2=$false},
    @{Row=31; Value="Plant pathogen"; Style2=$false},
    @{Row=32; Value="Plant pathogen"; Style2=$false},
    @{Row=33; Value="Plant pathogen"; Style2=$false},
    @{Row=34; Value="Plant pathogen"; Style2=$false},
    @{Row=35; Value="Plant pathogen"; Style2=$false},
    @{Row=36; Value="Plant pathogen"; Style2=$false},
    @{Row=37; Value="Saprotroph"; Style2=$false},
    @{Row=38; Value="Endophyte"; Style2=$false},
    @{Row=39; Value="Human pathogen"; Style2=$false},
    @{Row=40; Value="Carnivore"; Style2=$false},
    @{Row=41; Value="Plant pathogen"; Style2=$false},
    @{Row=42; Value="Plant pathogen"; Style2=$false},
    @{Row=43; Value="Plant pathogen"; Style2=$false},
    @{Row=44; Value="Human pathogen"; Style2=$false},
    @{Row=45; Value="Insect pathogen"; Style2=$false},
    @{Row=46; Value="Insect pathogen"; Style2=$false},
    @{Row=47; Value="Human pathogen"; Style2=$false},
    @{Row=48; Value="Saprotroph"; Style2=$true},
    @{Row=49; Value="Saprotroph"; Style2=$true},
    @{Row=50; Value="Saprotroph"; Style2=$true},
    @{Row=51; Value="Insect pathogen"; Style2=$false},
    @{Row=52; Value="Insect pathogen"; Style2=$false},
    @{Row=53; Value="Insect pathogen"; Style2=$false},
    @{Row=54; Value="Human pathogen"; Style2=$false},
    @{Row=55; Value="Human pathogen"; Style2=$false},
    @{Row=56; Value="Saprotroph"; Style2=$true},
    @{Row=57; Value="Plant pathogen"; Style2=$false},
    @{Row=58; Value="Plant pathogen"; Style2=$false},
    @{Row=59; Value="Saprotroph"; Style2=$true},
    @{Row=60; Value="Carnivore"; Style2=$false},
    @{Row=61; Value="Carnivore"; Style2=$false},
    @{Row=62; Value="Saprotroph"; Style2=$true},
    @{Row=63; Value="Carnivore"; Style2=$false},
    @{Row=64; Value="Plant pathogen"; Style2=$false},
    @{Row=65; Value="Plant pathogen"; Style2=$false},
    @{Row=66; Value="Plant pathogen"; Style2=$false},
    @{Row=67; Value="Saprotroph"; Style2=$true},
    @{Row=68; Value="Endophyte"; Style2=$false},
    @{Row=69; Value="Endophyte"; Style2=$false},
    @{Row=70; Value="Saprotroph"; Style2=$true},
    @{Row=71; Value="Plant pathogen"; Style2=$false},
    @{Row=72; Value="Plant pathogen"; Style2=$false},
    @{Row=73; Value="Saprotroph"; Style2=$true},
    @{Row=74; Value="Saprotroph"; Style2=$true},
    @{Row=75; Value="Saprotroph"; Style2=$true},
    @{Row=76; Value="Saprotroph"; Style2=$true},
    @{Row=77; Value="Saprotroph"; Style2=$true},
    @{Row=78; Value="Plant pathogen"; Style2=$false},
    @{Row=79; Value="Plant pathogen"; Style2=$false},
    @{Row=80; Value="Plant pathogen"; Style2=$false},
    @{Row=81; Value="Plant pathogen"; Style2=$false},
    @{Row=82; Value="Endophyte"; Style2=$false},
    @{Row=83; Value="Endophyte"; Style2=$false},
    @{Row=84; Value="Endophyte"; Style2=$false},
    @{Row=85; Value="Human pathogen"; Style2=$false},
    @{Row=86; Value="Plant pathogen"; Style2=$false},
    @{Row=87; Value="Saprotroph"; Style2=$false},
    @{Row=88; Value="Plant pathogen"; Style2=$false},
    @{Row=89; Value="Plant pathogen"; Style2=$false},
    @{Row=90; Value="Plant pathogen"; Style2=$false},
    @{Row=91; Value="Plant pathogen"; Style2=$false},
    @{Row=92; Value="Saprotroph"; Style2=$false},
    @{Row=93; Value="Saprotroph"; Style2=$false},
    @{Row=94; Value="Plant pathogen"; Style2=$false}
)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, 14)
    $cell.Value = $r.Value
    if ($r.Style2) {
        $cell.Font.Color = 0
    }
}

# --- Data cleanup: row 60 no longer lists "Plant" under Pathogen ------------
$ws.Range("I60").ClearContents() | Out-Null

# --- Restore view state (scrolled down, new cell selected) ------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A68").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("N94").Select() | Out-Null
